$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the sensitive password value from "Digitar Senha" row (C17) and
# the stray RA number from "Digitar RA" row (C13). Clearing the contents
# (rather than writing an empty string) drops the now-unused shared
# string "Harminda#2403" from the shared strings table, matching the
# target workbook.
$ws.Range("C13").ClearContents()
$ws.Range("C17").ClearContents()

# Update the active selection to C12.
[void]$ws.Range("C12").Select()
